$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 14.96897479739331; "C" = 8.266510995309281; "D" = 9.317208248789967; "F" = 36.31233528709617; "G" = 40.12048405863331; "H" = 16.83448487634694; "J" = 11.01257593831783; "M" = 17.96244524207853 }
    3 = @{ "B" = 14.46667559738459; "C" = 7.801617803030675; "D" = 9.310225729630393; "F" = 36.289288153946; "G" = 39.96665136336791; "H" = 16.87335582583014; "J" = 11.04035937685379; "M" = 17.81513519022682 }
    4 = @{ "B" = 14.15229384886916; "C" = 7.502268120912341; "D" = 9.306873983947732; "F" = 36.28654198038769; "G" = 39.88826225060501; "H" = 16.90135899576612; "J" = 11.05891815563326; "M" = 17.72756565199284 }
    5 = @{ "B" = 14.02288991811257; "C" = 7.376892821366991; "D" = 9.305744719364993; "F" = 36.28828528808188; "G" = 39.86037077594217; "H" = 16.91380642933412; "J" = 11.06685806005899; "M" = 17.69263686825482 }
    6 = @{ "B" = 14.00133100125933; "C" = 7.355873335982237; "D" = 9.305571533816375; "F" = 36.28874743498869; "G" = 39.85598443659587; "H" = 16.91593575592283; "J" = 11.06819924678885; "M" = 17.68688358416079 }
    7 = @{ "B" = 14.15055359874542; "C" = 7.500590823693167; "D" = 9.306857794574986; "F" = 36.28655391076212; "G" = 39.88786967520075; "H" = 16.90152267760562; "J" = 11.0590237093625; "M" = 17.72709148523725 }
    8 = @{ "B" = 14.79714127908995; "C" = 8.109159801432089; "D" = 9.31460732986625; "F" = 36.30201904710981; "G" = 40.06411830838726; "H" = 16.84702665903577; "J" = 11.02184439349466; "M" = 17.91107497550103 }
    9 = @{ "B" = 16.00940457210513; "C" = 9.188505318128692; "D" = 9.33716641645737; "F" = 36.42296654126203; "G" = 40.53629601505434; "H" = 16.77315647674327; "J" = 10.96083685303377; "M" = 18.2932292654396 }
    10 = @{ "B" = 16.85616472220004; "C" = 9.908118103799769; "D" = 9.358144642629613; "F" = 36.56707805424189; "G" = 40.95865935748063; "H" = 16.73922678941241; "J" = 10.92327269183548; "M" = 18.58496569587197 }
    11 = @{ "B" = 17.23008731050037; "C" = 10.2190052293421; "D" = 9.368625129090569; "F" = 36.6445799948486; "G" = 41.16668570274945; "H" = 16.72825139915164; "J" = 10.90776043148518; "M" = 18.71961312061453 }
    12 = @{ "B" = 17.36993584665775; "C" = 10.334329036844; "D" = 9.372726818170269; "F" = 36.67563633806266; "G" = 41.24769327909637; "H" = 16.7247395481868; "J" = 10.90211301401151; "M" = 18.77084019366104 }
    13 = @{ "B" = 17.33989655981606; "C" = 10.30959930039794; "D" = 9.371837562834058; "F" = 36.6688720158119; "G" = 41.23014845500069; "H" = 16.72546718828687; "J" = 10.90331920009545; "M" = 18.75979752265996 }
    14 = @{ "B" = 17.2416284367021; "C" = 10.22854130405968; "D" = 9.368959923670452; "F" = 36.6471008521399; "G" = 41.17330582365808; "H" = 16.72794954830107; "J" = 10.90729127013101; "M" = 18.72382303771801 }
    15 = @{ "B" = 17.18120524399695; "C" = 10.17857719718646; "D" = 9.367214545602559; "F" = 36.63398749037381; "G" = 41.13877715014267; "H" = 16.72955405563805; "J" = 10.90975380962893; "M" = 18.70181758696993 }
    16 = @{ "B" = 16.83148998138162; "C" = 9.887466396847147; "D" = 9.35747842531573; "F" = 36.56225272270022; "G" = 40.94537975247282; "H" = 16.74003405841492; "J" = 10.92431817630234; "M" = 18.57620175725233 }
    17 = @{ "B" = 16.61396873718504; "C" = 9.704634372287362; "D" = 9.351744433325143; "F" = 36.52129917430553; "G" = 40.83077124728808; "H" = 16.74760762259869; "J" = 10.93365665528284; "M" = 18.49960826497063 }
    18 = @{ "B" = 16.48780180972701; "C" = 9.597925619276127; "D" = 9.34853473877552; "F" = 36.49886883943979; "G" = 40.76635075130215; "H" = 16.75238331485768; "J" = 10.93917622174324; "M" = 18.45573821627696 }
    19 = @{ "B" = 16.44490686119794; "C" = 9.561531074022051; "D" = 9.347463219663549; "F" = 36.49146779172528; "G" = 40.74479808892016; "H" = 16.75407224726091; "J" = 10.94107052217274; "M" = 18.44091745296398 }
    20 = @{ "B" = 16.6372343436313; "C" = 9.724257702417685; "D" = 9.352345696167909; "F" = 36.5255423709082; "G" = 40.84281671105479; "H" = 16.74675795610472; "J" = 10.9326472067273; "M" = 18.5077429615516 }
    21 = @{ "B" = 17.27054053949351; "C" = 10.25241543653944; "D" = 9.36980156267056; "F" = 36.65344930526621; "G" = 41.18994174030879; "H" = 16.72720291056403; "J" = 10.90611842222576; "M" = 18.7343834424996 }
    22 = @{ "B" = 17.6742064324606; "C" = 10.58358346145337; "D" = 9.381984047075177; "F" = 36.7469932481436; "G" = 41.42979246197584; "H" = 16.71817915140824; "J" = 10.89010199310284; "M" = 18.88388102044535 }
    23 = @{ "B" = 17.45973709346439; "C" = 10.40812442318347; "D" = 9.375411825112003; "F" = 36.6961606062849; "G" = 41.30061038112275; "H" = 16.72265065640076; "J" = 10.89852928815895; "M" = 18.8039785296229 }
    24 = @{ "B" = 16.62671941389398; "C" = 9.71539096399491; "D" = 9.352073594571065; "F" = 36.52362054975008; "G" = 40.83736637574676; "H" = 16.74714077744992; "J" = 10.93310310885929; "M" = 18.50406474962931 }
    25 = @{ "B" = 15.68851059221284; "C" = 8.909176991014608; "D" = 9.330283791587485; "F" = 36.38053396384686; "G" = 40.39515912781748; "H" = 16.78958324534962; "J" = 10.9760666128696; "M" = 18.18776829137237 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
